$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F header ---
$ws.Cells.Item(1, 6).Value = "Trening"
# Reuse the existing header style (s="1") instead of creating a new one.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# --- Convert column A (rows 2-7) from text timestamps to real datetime serials ---
# and tag the existing rows as belonging to "Duza Gra" (big game) segment.
$existingTimestamps = @{
    2 = 45672.47361111111
    3 = 45672.48819444444
    4 = 45672.49375
    5 = 45672.46458333333
    6 = 45672.47013888889
    7 = 45672.47361111111
}

# Priming step: touch the first date cell with a lowercase format code first
# (registers numFmtId 164), then switch it to the uppercase variant that is
# actually used throughout the sheet (registers numFmtId 165). Every other
# date cell below is set directly to the uppercase format so they all share
# the same style index.
$primeCell = $ws.Cells.Item(2, 1)
$primeCell.Value = $existingTimestamps[2]
$primeCell.NumberFormat = "yyyy-mm-dd h:mm:ss"
$primeCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 6).Value = "Duża Gra"

foreach ($r in 3..7) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $existingTimestamps[$r]
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 6).Value = "Duża Gra"
}

# --- Append new rows 8-13 for the "Mala Gra" (small game) segment ---
$newRows = @(
    @{ Row = 8;  A = 45672.46388888889; B = 279.1; C = 10.87; D = 0.9338212012857143;  E = "10-15" }
    @{ Row = 9;  A = 45672.46458333333; B = 340.6; C = 11.29; D = 3.156308037857143;   E = "10-15" }
    @{ Row = 10; A = 45672.46527777778; B = 392.2; C = 12.3;  D = 3.135422025571429;   E = "10-15" }
    @{ Row = 11; A = 45672.46458333333; B = 340.4; C = 8.99;  D = 2.996093545714285;   E = "5-10" }
    @{ Row = 12; A = 45672.46527777778; B = 388.1; C = 6.02;  D = 1.776477507142857;   E = "5-10" }
    @{ Row = 13; A = 45672.46527777778; B = 391.9; C = 8.73;  D = 2.766690067142857;   E = "5-10" }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $item.A
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = "Mała Gra"
}
